$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text could be misread as a number by Excel's COM Value setter;
# force them to Text format first, write the value, then restore the default style
# (so the saved file has no stray numFmt on these cells, matching the original styling).
function Set-TextValue($addr, $text) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).Style = "Normal"
}

$ws.Range('D2').Value = '61.465.78'
$ws.Range('E2').Value = '  +0.58%  '
$ws.Range('D3').Value = '2.933.82'
$ws.Range('E3').Value = '  +0.14%  '
$ws.Range('E4').Value = '  -0.02%  '
Set-TextValue 'D5' '597.01'
$ws.Range('E5').Value = '  +0.93%  '
Set-TextValue 'D6' '145.32'
$ws.Range('E6').Value = '  -0.63%  '
$ws.Range('E7').Value = '  +0.02%  '
Set-TextValue 'D8' '0.502'
$ws.Range('E8').Value = '  -0.71%  '
Set-TextValue 'D9' '7.01'
$ws.Range('E9').Value = '  +1.94%  '
$ws.Range('E10').Value = '  -1.65%  '
Set-TextValue 'D11' '0.439'
$ws.Range('E11').Value = '  -0.50%  '
$ws.Range('E12').Value = '  -0.85%  '
Set-TextValue 'D13' '33.59'
$ws.Range('E13').Value = '  -0.58%  '
$ws.Range('E14').Value = '  +0.50%  '
$ws.Range('D15').Value = '3.420.72'
$ws.Range('E15').Value = '  +0.21%  '
$ws.Range('D16').Value = '61.439.97'
$ws.Range('E16').Value = '  +0.61%  '
$ws.Range('D17').Value = '2.934.64'
$ws.Range('E17').Value = '  +0.25%  '
$ws.Range('E18').Value = '  -0.09%  '
Set-TextValue 'D19' '432.01'
$ws.Range('E19').Value = '  +0.09%  '
Set-TextValue 'D20' '13.48'
$ws.Range('E20').Value = '  +0.26%  '
Set-TextValue 'D21' '0.676'
$ws.Range('E21').Value = '  -1.09%  '
Set-TextValue 'D22' '7.10'
$ws.Range('E22').Value = '  -0.03%  '
Set-TextValue 'D23' '82.01'
$ws.Range('E23').Value = '  +0.77%  '
$ws.Range('E24').Value = '  -1.07%  '
Set-TextValue 'D25' '2.19'
$ws.Range('E25').Value = '  -1.62%  '
Set-TextValue 'D26' '11.78'
$ws.Range('E26').Value = '  -2.29%  '
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('E28').Value = '  -4.16%  '
$ws.Range('E29').Value = '  -0.21%  '
Set-TextValue 'D30' '6.94'
$ws.Range('E30').Value = '  -2.20%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D31' '26.66'
$ws.Range('E31').Value = '  +0.60%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D32' '0.110'
$ws.Range('E32').Value = '  +1.24%  '
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('D34').Value = '0.0₃0885'
$ws.Range('E34').Value = '  +2.75%  '
$ws.Range('E35').Value = '  +0.04%  '
Set-TextValue 'D36' '5.63'
$ws.Range('E36').Value = '  +0.18%  '
Set-TextValue 'D37' '2.99'
$ws.Range('E37').Value = '  -2.61%  '
$ws.Range('E38').Value = '  +0.15%  '
$ws.Range('E39').Value = '  -1.63%  '
Set-TextValue 'D40' '8.58'
$ws.Range('E40').Value = '  -0.13%  '
Set-TextValue 'D41' '42.33'
$ws.Range('E41').Value = '  +7.28%  '
Set-TextValue 'D42' '0.282'
$ws.Range('E42').Value = '  -1.74%  '
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('D44').Value = '2.697.00'
$ws.Range('E44').Value = '  -0.64%  '
Set-TextValue 'D45' '134.22'
$ws.Range('E45').Value = '  +1.86%  '
Set-TextValue 'D46' '362.04'
$ws.Range('E46').Value = '  -3.62%  '
$ws.Range('E47').Value = '  +0.03%  '
Set-TextValue 'D48' '23.70'
$ws.Range('E48').Value = '  -2.17%  '
$ws.Range('E49').Value = '  -1.37%  '
Set-TextValue 'D50' '2.00'
$ws.Range('E50').Value = '  -1.71%  '
$ws.Range('E51').Value = '  -1.43%  '
